# Update the "Example" column values (E, H, K) for rows 5-9 to reflect
# the refreshed descriptor examples.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


$ws.Range("E5").Value = "Molecular weight"
$ws.Range("E6").Value = "Molecular weight"
$ws.Range("E7").Value = "Molecular weight"
$ws.Range("E8").Value = "Molecular weight`nNumber of smallest rings`nNumber of acidic groups"
$ws.Range("E9").Value = "Molecular weight`nInChIKey`nSMARTS filter`nSpectrophores"
$ws.Range("H5").Value = "Valence molecular connectivity Chi index for path order 0`nWiener Index`nThe information content with order 0 proposed by Basak`nSecond kappa shape index`nHighest eigenvaluen.1 of Burden matrix/weighted by atomic masses`nSum of E-State of atom type: sLi`nMoreau-Broto weight autocorrelation by mass`nMolar refractivity`nMost positive charge on H atoms `nMOE Charge VSA Descriptor 1"
$ws.Range("H6").Value = "Valence molecular connectivity Chi index for path order 0`nMolar refractivity`nMOE Charge VSA Descriptor 1`nSecond kappa shape index`nSum of E-State of atom type: sLi`nWiener Index"
$ws.Range("H7").Value = "Wiener Index`nMolecular walk count of order 1`nValence molecular connectivity Chi index for path order 0`nTotal information index on atomic composition`nBalaban-like index from adjacency matrix`nMoreau-Broto weight autocorrelation by mass`nHighest eigenvaluen.1 of Burden matrix/weighted by atomic masses`nMOE Charge VSA Descriptor 1`nETA core count`nLeading eigenvalue from edge adjacency matrix`nFrequency of P – P at topological distance 6`nSum of E-State of atom type: sLi`nHydrogen attached to heteroatom`nSHED Acceptor-Negative`nMost positive charge on H atoms `nMolar refractivity"
$ws.Range("H8").Value = "Moreau-Broto weight autocorrelation by mass`nHighest eigenvaluen.1 of Burden matrix/weighted by atomic masses`nValence molecular connectivity Chi index for path order 0`nSecond kappa shape index`nSum of E-State of atom type: sLi`nMolar refractivity`nIonization potential of a molecule`nWiener Index"
$ws.Range("K5").Value = "Radius of gyration`nPartial negative surface area`n3D-RDF - signal 01 / unweighted`nCalculate Wlambda1`n3D-MoRSE - signal 01 / weighted by atomic charge"
$ws.Range("K6").Value = "Radius of gyration"
$ws.Range("K7").Value = "Radius of gyration`nWiener-like index from geometrical matrix`n3D Topological distance based descriptors – lag 1 unweighted`nRadial Distribution Function – 155 / weighted by I-state`n3D-MoRSE - signal 01 / weighted by atomic charge`nCalculate Wlambda1`nR total index / weighted by I-state`nMolecular profile no. 1`nNumber of aromatic aldehydes`nSum of geometrical distances between S and F`nQuantitative Estimation of Drug-likeness`nCATS3D Lipophilic-Lipophilic BIN 19`nWHALES Remoteness`nMolecular distance edge between all primary carbons`ns1_phSize normalized by the heavy atoms"
$ws.Range("K8").Value = "Partial negative surface area`nCalculate Wlambda1`nMolecular distance edge between all primary carbons`nRadius of gyration"

# Move the active selection to K7 (was K8).
$ws.Range("K7").Select()
